$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs)
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.022418
$ws.Range("H2").Value = 0.06725399999999999
$ws.Range("M2").Value = 15.24491733333333
$ws.Range("N2").Value = 45.73475199999999
$ws.Range("O2").Value = 0.4831257321597052
$ws.Range("P2").Value = 0.4831257321597052
$ws.Range("Q2").Value = 0.3417605567786666
$ws.Range("R2").Value = 3.075845011007999
$ws.Range("S2").Value = 0.4831257321597052
$ws.Range("T2").Value = 0.4831257321597052

# Row 3 (FAPs)
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.022418
$ws.Range("H3").Value = 0.06725399999999999
$ws.Range("O3").Value = 0.327710667227878
$ws.Range("P3").Value = 0.327710667227878
$ws.Range("Q3").Value = 0.231820771776
$ws.Range("R3").Value = 2.086386945984
$ws.Range("S3").Value = 0.327710667227878
$ws.Range("T3").Value = 0.327710667227878

# Row 4 (MuSCs)
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.022418
$ws.Range("H4").Value = 0.06725399999999999
$ws.Range("M4").Value = 5.969012333333333
$ws.Range("N4").Value = 17.907037
$ws.Range("O4").Value = 0.1891636006124168
$ws.Range("P4").Value = 0.1891636006124168
$ws.Range("Q4").Value = 0.1338133184886666
$ws.Range("R4").Value = 1.204319866398
$ws.Range("S4").Value = 0.1891636006124168
$ws.Range("T4").Value = 0.1891636006124168
